# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Ciruela" variety "Angeleno" at
# Macroferia Regional de Talca, pushing the existing rows (old 96..112) down
# to become rows 98..114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 96 (old row96 + everything
# below shifts down by 2 rows total, one insert at a time).
$ws.Rows.Item(96).Insert()
$ws.Rows.Item(96).Insert()

# New row 96: Angeleno / Primera
$ws.Range("A96").Value = 5
$ws.Range("B96").Value = "Macroferia Regional de Talca"
$ws.Range("C96").Value = "Maule"
$ws.Range("D96").Value = 44644
$ws.Range("E96").Value = 7
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100103
$ws.Range("H96").Value = "Frutos de hueso (carozo)"
$ws.Range("I96").Value = 100103002
$ws.Range("J96").Value = "Ciruela"
$ws.Range("K96").Value = "Angeleno"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 200
$ws.Range("N96").Value = 8000
$ws.Range("O96").Value = 8000
$ws.Range("P96").Value = 8000
$ws.Range("Q96").Value = "$/bandeja 18 kilos granel"
$ws.Range("R96").Value = "Provincia de Curicó"
$ws.Range("S96").Value = 444
$ws.Range("T96").Value = 18

# New row 97: Angeleno / Segunda
$ws.Range("A97").Value = 5
$ws.Range("B97").Value = "Macroferia Regional de Talca"
$ws.Range("C97").Value = "Maule"
$ws.Range("D97").Value = 44644
$ws.Range("E97").Value = 7
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100103
$ws.Range("H97").Value = "Frutos de hueso (carozo)"
$ws.Range("I97").Value = 100103002
$ws.Range("J97").Value = "Ciruela"
$ws.Range("K97").Value = "Angeleno"
$ws.Range("L97").Value = "Segunda"
$ws.Range("M97").Value = 100
$ws.Range("N97").Value = 6000
$ws.Range("O97").Value = 6000
$ws.Range("P97").Value = 6000
$ws.Range("Q97").Value = "$/bandeja 18 kilos granel"
$ws.Range("R97").Value = "Provincia de Curicó"
$ws.Range("S97").Value = 333
$ws.Range("T97").Value = 18
